# Append the 2017-08-04 (serial 42951) temperature/tote readings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 42951
$toteValues = @(9,7,6,5,5,9,3,3,4,7,2,8,8,4,2,1,6,8,7,9,3,4,2,6)

$startRow = 74
for ($i = 0; $i -lt $toteValues.Length; $i++) {
    $row = $startRow + $i
    $bucket = $i + 1

    # Copy the previous row's formatting (date style in column A) down before
    # writing values, so the new cells inherit style index 2 like the rest
    # of the date column instead of Excel minting a fresh numFmt xf.
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $newDate
    $ws.Cells.Item($row, 2).Value = $bucket
    $ws.Cells.Item($row, 3).Value = "N/A"
    $ws.Cells.Item($row, 4).Value = $toteValues[$i]
}

# Match the saved view state from the edit (scrolled near the bottom, D93 selected,
# zoomed to 102%). The engine doesn't persist window scroll position
# (topLeftCell) to the saved sheetView, but zoom and selection do round-trip.
$ws.Range("D93").Select()
$excel.ActiveWindow.ScrollRow = 71
$excel.ActiveWindow.Zoom = 102
